$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New bet rows (122-126) -- fill in the previously-blank placeholder rows.
# Column A needs the date number-format (style used by A117:A121) copied in,
# everything else already carries the right style on the blank rows so only
# values/formulas need to be written.
# ---------------------------------------------------------------------------

$ws.Range("A121").Copy()
$ws.Range("A122:A126").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Rows are entered with every game marked "Pending" first (matching the
# order new shared strings were minted: event name, bet name, then the
# freshly-typed "Pending" status for row 122, reused by the later rows).
# Finished games are then flipped to their final result afterwards, which
# reuses the pre-existing "Win"/"Loss" shared strings, so row order below
# mirrors the original authoring sequence.

# Row 122 - Troy vs UTSA (Moneyline, UTSA)
$ws.Range("A122").Value = 44909
$ws.Range("B122").Value = "Troy vs UTSA"
$ws.Range("C122").Value = "Moneyline"
$ws.Range("D122").Value = "UTSA"
$ws.Range("E122").Value = 105
$ws.Range("F122").Value = 105
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 2.05
$ws.Range("I122").Value = "Pending"

# Row 123 - Stony Brook vs Wagner (Moneyline, Stony Brook)
$ws.Range("A123").Value = 44910
$ws.Range("B123").Value = "Stony Brook vs Wagner"
$ws.Range("C123").Value = "Moneyline"
$ws.Range("D123").Value = "Stony Brook"
$ws.Range("E123").Value = 230
$ws.Range("F123").Value = 240
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 3.3
$ws.Range("I123").Value = "Pending"

# Row 124 - Dallas Stars vs Carolina Hurricanes (Point Spread, Dallas Stars +1.5, Pending)
$ws.Range("A124").Value = 44911
$ws.Range("B124").Value = "Dallas Stars vs Carolina Hurricanes"
$ws.Range("C124").Value = "Point Spread"
$ws.Range("D124").Value = "Dallas Stars +1.5"
$ws.Range("E124").Value = -175
$ws.Range("F124").Value = -175
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 1.57
$ws.Range("I124").Value = "Pending"

# Row 125 - Clemson vs Tennessee (Moneyline, Tennessee, Pending)
$ws.Range("A125").Value = 44911
$ws.Range("B125").Value = "Clemson vs Tennessee"
$ws.Range("C125").Value = "Moneyline"
$ws.Range("D125").Value = "Tennessee"
$ws.Range("E125").Value = 215
$ws.Range("F125").Value = 215
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 3.15
$ws.Range("I125").Value = "Pending"

# Row 126 - Florida vs Oregon State (Moneyline, Florida, Pending)
$ws.Range("A126").Value = 44912
$ws.Range("B126").Value = "Florida vs Oregon State"
$ws.Range("C126").Value = "Moneyline"
$ws.Range("D126").Value = "Florida"
$ws.Range("E126").Value = 260
$ws.Range("F126").Value = 260
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 3.6
$ws.Range("I126").Value = "Pending"

# Troy vs UTSA and Stony Brook vs Wagner have since finished -- flip their
# results to Loss (reuses the pre-existing "Loss" shared string, index 5).
$ws.Range("I122").Value = "Loss"
$ws.Range("I123").Value = "Loss"

# ---------------------------------------------------------------------------
# Extend the shared formulas down through row 126 (they previously stopped
# at row 121).
# ---------------------------------------------------------------------------
$ws.Range("J122:J126").Formula = '=IF(I122="Pending", 0,IF(I122="Win",H122-G122,-1*G122))'
$ws.Range("K122:K126").Formula = '=K121+J122'

# ---------------------------------------------------------------------------
# View state: active cell / scroll position, matching the author's saved
# selection after adding the rows above.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 104
$win.ScrollColumn = 1
$ws.Range("L124").Select()
